# Edit script: update date + table cell contents (regenerated practice problems)
$d = $word.ActiveDocument

# Update the date line (paragraph 1, before the table)
$d.Paragraphs(1).Range.Text = "2024-06-24 Monday"

$t = $d.Tables(1)

# Table row 1 (filled data row)
$t.Cell(1, 1).Range.Text = "88÷8=11, 0"
$t.Cell(1, 2).Range.Text = "80÷8=10, 0"
$t.Cell(1, 3).Range.Text = "61÷2=30, 1"
$t.Cell(1, 4).Range.Text = "51÷2=25, 1"
$t.Cell(1, 5).Range.Text = "57÷4=14, 1"

# Table row 5 (filled data row)
$t.Cell(5, 1).Range.Text = "62÷6=10, 2"
$t.Cell(5, 2).Range.Text = "50÷2=25, 0"
$t.Cell(5, 3).Range.Text = "58÷6=9, 4"
$t.Cell(5, 4).Range.Text = "68÷8=8, 4"
$t.Cell(5, 5).Range.Text = "90÷8=11, 2"

# Table row 9 (filled data row)
$t.Cell(9, 1).Range.Text = "17÷3=5, 2"
$t.Cell(9, 2).Range.Text = "40÷6=6, 4"
$t.Cell(9, 3).Range.Text = "70÷9=7, 7"
$t.Cell(9, 4).Range.Text = "97÷8=12, 1"
$t.Cell(9, 5).Range.Text = "84÷5=16, 4"

# Table row 13 (filled data row)
$t.Cell(13, 1).Range.Text = "79÷9=8, 7"
$t.Cell(13, 2).Range.Text = "53÷8=6, 5"
$t.Cell(13, 3).Range.Text = "72÷9=8, 0"
$t.Cell(13, 4).Range.Text = "87÷7=12, 3"
$t.Cell(13, 5).Range.Text = "58÷9=6, 4"

# Table row 17 (filled data row)
$t.Cell(17, 1).Range.Text = "53÷4=13, 1"
$t.Cell(17, 2).Range.Text = "69÷8=8, 5"
$t.Cell(17, 3).Range.Text = "51÷7=7, 2"
$t.Cell(17, 4).Range.Text = "34÷7=4, 6"
$t.Cell(17, 5).Range.Text = "52÷6=8, 4"
